# Update specific data values in the first worksheet.
# These correspond to result_data_KNN.xlsx imputed values being refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.523
$ws.Range("E4").Value = 13.265
$ws.Range("B9").Value = 6.764
$ws.Range("E10").Value = 12.554
$ws.Range("B18").Value = 5.972
$ws.Range("B20").Value = 6.37
$ws.Range("D21").Value = -7.805
